# Update the Cgn -> F11r ligand-receptor pair sheet with freshly recomputed
# TPM-based NATMI statistics. The old table only had a single "Sending
# cluster" (MuSCs, rows 2-5). The refreshed run adds "ECs" as a second
# sending cluster and recomputes every numeric column for both senders
# against each of the four target clusters (ECs, FAPs, MuSCs, Resolving-Mac).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each row: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric NATMI metric columns (E..T).
$rows = @(
    @("ECs",   "Cgn", "F11r", "ECs",           1, 0.3333333333333333, 0.02557233333333333, 0.076717,            0.6935308900902204, 0.6935308900902204, 3, 1,                  40.91514966666667, 122.745449, 0.8529192913871414,  0.8529192913871415,  1.046295845659222,   9.416662610932999, 0.5915258753308442,  0.5915258753308443),
    @("ECs",   "Cgn", "F11r", "FAPs",          1, 0.3333333333333333, 0.02557233333333333, 0.076717,            0.6935308900902204, 0.6935308900902204, 2, 0.6666666666666666, 0.165216,           0.495648,   0.00344410114086962,  0.003444101140869621, 0.004224958624,      0.038024627616,    0.002388590529788051, 0.002388590529788052),
    @("ECs",   "Cgn", "F11r", "MuSCs",         1, 0.3333333333333333, 0.02557233333333333, 0.076717,            0.6935308900902204, 0.6935308900902204, 2, 0.6666666666666666, 0.4441646666666667, 1.332494,   0.009259079236881667, 0.009259079236881667, 0.01135832691088889, 0.102224942198,    0.006421457464570421, 0.006421457464570421),
    @("ECs",   "Cgn", "F11r", "Resolving-Mac", 1, 0.3333333333333333, 0.02557233333333333, 0.076717,            0.6935308900902204, 0.6935308900902204, 3, 1,                  6.446186333333333, 19.338559,  0.1343775282351073,  0.1343775282351073,  0.1648440256447778,  1.483596230803,    0.09319496676501768, 0.0931949667650177),
    @("MuSCs", "Cgn", "F11r", "ECs",           1, 0.3333333333333333, 0.01130033333333333, 0.033901,            0.3064691099097797, 0.3064691099097796, 3, 1,                  40.91514966666667, 122.745449, 0.8529192913871414,  0.8529192913871415,  0.4623548296165557,  4.161193466549,    0.2613934160562972, 0.2613934160562972),
    @("MuSCs", "Cgn", "F11r", "FAPs",          1, 0.3333333333333333, 0.01130033333333333, 0.033901,            0.3064691099097797, 0.3064691099097796, 2, 0.6666666666666666, 0.165216,           0.495648,   0.00344410114086962,  0.003444101140869621, 0.001866995872,      0.016802962848,    0.001055510611081569, 0.001055510611081569),
    @("MuSCs", "Cgn", "F11r", "MuSCs",         1, 0.3333333333333333, 0.01130033333333333, 0.033901,            0.3064691099097797, 0.3064691099097796, 2, 0.6666666666666666, 0.4441646666666667, 1.332494,   0.009259079236881667, 0.009259079236881667, 0.005019208788222223, 0.045172879094,    0.002837621772311246, 0.002837621772311246),
    @("MuSCs", "Cgn", "F11r", "Resolving-Mac", 1, 0.3333333333333333, 0.01130033333333333, 0.033901,            0.3064691099097797, 0.3064691099097796, 3, 1,                  6.446186333333333, 19.338559,  0.1343775282351073,  0.1343775282351073,  0.07284405429544445, 0.655596488659,    0.04118256147008961, 0.04118256147008962)
)

$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $addr = "$($columns[$i])$r"
        $ws.Range($addr).Value = $row[$i]
    }
    $r++
}
